# Correção nos dados e inicio da analise PNAD 2009
#
# The original sheet had two "section header" rows (row 5 = "situação do
# domicílio", row 8 = "grandes regiões e unidades da federação") that only
# carried a label in column A with no data in B:G. The corrected data
# removes those two header-only rows entirely (their data lives in the
# rows immediately below, which already hold real numbers), so every row
# below them shifts up. The top "unnamed: 1_level_1" column header above
# "própria residência ou de terceiros" is also fixed to read "total" (it
# mirrors the "total" header already used in B1), matching the corrected
# header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two label-only rows; this shifts everything below them up,
# carrying the existing (correct) numeric data into the newly vacated
# row numbers, and the workbook engine automatically drops the
# now-unreferenced shared strings ("situação do domicílio" and
# "grandes regiões e unidades da federação") from sharedStrings.xml.
$ws.Rows(5).Delete()
$ws.Rows(7).Delete()

# Fix the mislabeled column header in row 2 (was "unnamed: 1_level_1").
$ws.Range("B2").Value = "total"
